$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (C) column from 45182 to 45184 for all data rows (2..439)
for ($i = 2; $i -le 439; $i++) {
    $ws.Cells.Item($i, 3).Value = 45184
}

# 2. Row 439 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(439).RowHeight = 15

# 3. Append four new data rows (440-443)

# --- Row 440 ---
$r = 440
$ws.Rows.Item($r).RowHeight = 15
$ws.Cells.Item($r, 1).Value = "A 42574-2023"
$ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($r, 2).Value = 45181
$ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($r, 3).Value = 45184
$ws.Cells.Item($r, 4).Value = "KALMAR LÄN"
$ws.Cells.Item($r, 5).Value = "HÖGSBY"
$ws.Cells.Item($r, 7).Value = 1.3
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0
$ws.Cells.Item($r, 11).Value = 0
$ws.Cells.Item($r, 12).Value = 0
$ws.Cells.Item($r, 13).Value = 0
$ws.Cells.Item($r, 14).Value = 0
$ws.Cells.Item($r, 15).Value = 0
$ws.Cells.Item($r, 16).Value = 0
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).WrapText = $true

# --- Row 441 ---
$r = 441
$ws.Rows.Item($r).RowHeight = 15
$ws.Cells.Item($r, 1).Value = "A 42668-2023"
$ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($r, 2).Value = 45181
$ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($r, 3).Value = 45184
$ws.Cells.Item($r, 4).Value = "KALMAR LÄN"
$ws.Cells.Item($r, 5).Value = "HÖGSBY"
$ws.Cells.Item($r, 6).Value = "Sveaskog"
$ws.Cells.Item($r, 7).Value = 13.4
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0
$ws.Cells.Item($r, 11).Value = 0
$ws.Cells.Item($r, 12).Value = 0
$ws.Cells.Item($r, 13).Value = 0
$ws.Cells.Item($r, 14).Value = 0
$ws.Cells.Item($r, 15).Value = 0
$ws.Cells.Item($r, 16).Value = 0
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).WrapText = $true

# --- Row 442 ---
$r = 442
$ws.Rows.Item($r).RowHeight = 15
$ws.Cells.Item($r, 1).Value = "A 42569-2023"
$ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($r, 2).Value = 45181
$ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($r, 3).Value = 45184
$ws.Cells.Item($r, 4).Value = "KALMAR LÄN"
$ws.Cells.Item($r, 5).Value = "HÖGSBY"
$ws.Cells.Item($r, 7).Value = 2.3
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0
$ws.Cells.Item($r, 11).Value = 0
$ws.Cells.Item($r, 12).Value = 0
$ws.Cells.Item($r, 13).Value = 0
$ws.Cells.Item($r, 14).Value = 0
$ws.Cells.Item($r, 15).Value = 0
$ws.Cells.Item($r, 16).Value = 0
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).WrapText = $true

# --- Row 443 (no explicit row height, matches diff) ---
$r = 443
$ws.Cells.Item($r, 1).Value = "A 43213-2023"
$ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($r, 2).Value = 45183
$ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($r, 3).Value = 45184
$ws.Cells.Item($r, 4).Value = "KALMAR LÄN"
$ws.Cells.Item($r, 5).Value = "HÖGSBY"
$ws.Cells.Item($r, 7).Value = 1.8
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0
$ws.Cells.Item($r, 11).Value = 0
$ws.Cells.Item($r, 12).Value = 0
$ws.Cells.Item($r, 13).Value = 0
$ws.Cells.Item($r, 14).Value = 0
$ws.Cells.Item($r, 15).Value = 0
$ws.Cells.Item($r, 16).Value = 0
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).WrapText = $true

Write-Host "Edit complete"
